$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Hend Mahmoud, Dr. Alshimaa Atef, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali'
$ws.Range("G3").Value = 'Administrator, Dr. Alshimaa Atef, Dr. Gehan Adel, Dr. Manar Montaser'
$ws.Range("G4").Value = 'Dr. Hanan Ragab, Dr. Majorelle Magdy, Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki, Dr. Nourhan Mahmoud, Dr. Asmaa Reda, Dr. Heba Mahmoud Ali'
$ws.Range("G6").Value = 'Dr. Safa Hany, Dr. Sara Nabil'
$ws.Range("G10").Value = 'Dr. Amira Ibrahim, Dr. Basma Hamed'
$ws.Range("G12").Value = 'Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range("G17").Value = 'Dr. Walaa Ghanima, Dr. Enas Omran, Dr. Marian Samir'
$ws.Range("G18").Value = 'Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida'
$ws.Range("G19").Value = 'Dr. Neveen Nashaat, Dr. Wafaa Ebida, Dr. Marina Sorial, Dr. Yasmin, Dr. Eman Samir Gabry'
$ws.Range("G20").Value = 'Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Yasmin, Dr. Remon, Dr. Monica, Dr. Marina Atef, Dr. Nardine'
$ws.Range("G21").Value = 'Dr. Hend Mahmoud, Dr. Alshimaa Atef, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali'
$ws.Range("G22").Value = 'Administrator, Dr. Alshimaa Atef, Dr. Gehan Adel, Dr. Manar Montaser'
$ws.Range("G23").Value = 'Dr. Hanan Ragab, Dr. Majorelle Magdy, Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki, Dr. Nourhan Mahmoud, Dr. Asmaa Reda, Dr. Heba Mahmoud Ali'
$ws.Range("G24").Value = 'Dr. Lamiaa Ossama, Dr. Amera Ahmad Saad, Dr. Abeer Ragab, Dr. Nada Mohammad, Dr. Fatma Elhady'
$ws.Range("G28").Value = 'Dr. Eman M. Abo-Sakaya, Dr. Dina Adel, Dr. Marwa Mustafa, Dr. Esraa Mostafa, Dr. Yasmeena Fattoh, Dr. Nourhan Osama, Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Sarah Abdelmohsen, Dr. Arwa Al-Sayed'
$ws.Range("G31").Value = 'Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range("G32").Value = 'Dr. Nouran Mahmoud, Menna tuâ€™Allah Gamil'
$ws.Range("G36").Value = 'Dr. Walaa Ghanima, Dr. Enas Omran, Dr. Marian Samir'
$ws.Range("G37").Value = 'Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida'
$ws.Range("G38").Value = 'Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Yasmin, Dr. Remon, Dr. Monica, Dr. Marina Atef, Dr. Nardine'
$ws.Range("G39").Value = 'Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Yasmin, Dr. Remon, Dr. Monica, Dr. Marina Atef, Dr. Nardine'
$ws.Range("G40").Value = 'Dr. Hend Mahmoud, Dr. Alshimaa Atef, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali'
$ws.Range("G41").Value = 'Dr. Hanan Ragab, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Hend Mahmoud, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki'
$ws.Range("G42").Value = 'Dr. Alshimaa Atef, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki'
$ws.Range("G43").Value = 'Dr. Lamiaa Ossama, Dr. Kerelos Zareef, Dr. Amera Ahmad Saad, Dr. Abeer Ragab, Dr. Menna tu''Alllah Mohammad, Dr. Nada Mohammad, Dr. Fatma Elhady'
$ws.Range("G44").Value = 'Dr. Safa Hany, Dr. Sara Nabil'
$ws.Range("G47").Value = 'Dr. Merna Said, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Maryam Ahmad, Dr. Amira Ibrahim, Dr. Arwa Al-Sayed'
$ws.Range("G48").Value = 'Dr. Merna Said, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Fatma Shoukry, Dr. Maryam Ahmad, Dr. Sarah Abdelmohsen, Dr. Amany Raafat'
$ws.Range("G49").Value = 'Dr. Mohammad Safwat, Dr. Mariam Toma Gerges'
$ws.Range("G50").Value = 'Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range("G51").Value = 'Dr. Nouran Mahmoud, Menna tuâ€™Allah Gamil'
$ws.Range("G54").Value = 'Dr. Amr Saeed, Dr. Afaf Abdallah'
$ws.Range("G56").Value = 'Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida'
$ws.Range("G57").Value = 'Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Yasmin, Dr. Remon, Dr. Monica, Dr. Marina Atef, Dr. Nardine'
$ws.Range("G58").Value = 'Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Yasmin, Dr. Remon, Dr. Monica, Dr. Marina Atef, Dr. Nardine'
$ws.Range("G59").Value = 'Dr. Amira Sobhy, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad, Dr. Nesma, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Asmaa Reda'
$ws.Range("G60").Value = 'Dr. Hanan Ragab, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Hend Mahmoud, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki'
$ws.Range("G61").Value = 'Dr. Amira Sobhy, Dr. Nahla Nagiub, Dr. Majorelle Magdy, Dr. Shimaa Ahmad Mekki, Dr. Asmaa Reda'
$ws.Range("G63").Value = 'Dr. Amal Awwad, Dr. Safa Hany, Dr. Aya Saeed'
$ws.Range("G66").Value = 'Dr. Eman M. Abo-Sakaya, Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Amira Ibrahim, Dr. Madeha Saeed, Dr. Marina Youhanna'
$ws.Range("G75").Value = 'Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida'
$ws.Range("G76").Value = 'Dr. Neveen Nashaat, Dr. Wafaa Ebida, Dr. Marina Sorial, Dr. Yasmin, Dr. Eman Samir Gabry'
$ws.Range("G77").Value = 'Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Yasmin, Dr. Remon, Dr. Monica, Dr. Marina Atef, Dr. Nardine'
$ws.Range("G78").Value = 'Dr. Amira Sobhy, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad, Dr. Nesma, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Asmaa Reda'
$ws.Range("G79").Value = 'Dr. Hanan Ragab, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Hend Mahmoud, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki'
$ws.Range("G80").Value = 'Dr. Amira Sobhy, Dr. Nahla Nagiub, Dr. Majorelle Magdy, Dr. Shimaa Ahmad Mekki, Dr. Asmaa Reda'
$ws.Range("G81").Value = 'Dr. Lamiaa Ossama, Dr. Amera Ahmad Saad, Dr. Abeer Ragab, Dr. Nada Mohammad, Dr. Fatma Elhady'
$ws.Range("G83").Value = 'Dr. Amal Awwad, Dr. Safa Hany, Dr. Aya Saeed'
$ws.Range("G85").Value = 'Dr. Eman M. Abo-Sakaya, Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Amira Ibrahim, Dr. Madeha Saeed, Dr. Marina Youhanna'
$ws.Range("G86").Value = 'Dr. Merna Said, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Fatma Shoukry, Dr. Maryam Ahmad, Dr. Sarah Abdelmohsen, Dr. Amany Raafat'
$ws.Range("G88").Value = 'Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range("G94").Value = 'Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida'
$ws.Range("G95").Value = 'Dr. Neveen Nashaat, Dr. Wafaa Ebida, Dr. Marina Sorial, Dr. Yasmin, Dr. Eman Samir Gabry'
$ws.Range("G96").Value = 'Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Yasmin, Dr. Remon, Dr. Monica, Dr. Marina Atef, Dr. Nardine'
$ws.Range("G97").Value = 'Dr. Amira Sobhy, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad, Dr. Nesma, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Asmaa Reda'
$ws.Range("G98").Value = 'Dr. Hanan Ragab, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Hend Mahmoud, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki'
$ws.Range("G99").Value = 'Dr. Alshimaa Atef, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki'
$ws.Range("G100").Value = 'Dr. Lamiaa Ossama, Dr. Kerelos Zareef, Dr. Amera Ahmad Saad, Dr. Abeer Ragab, Dr. Menna tu''Alllah Mohammad, Dr. Nada Mohammad, Dr. Fatma Elhady'
$ws.Range("G101").Value = 'Dr. Amal Awwad, Dr. Safa Hany, Dr. Aya Saeed'
$ws.Range("G104").Value = 'Dr. Merna Said, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Maryam Ahmad, Dr. Amira Ibrahim, Dr. Arwa Al-Sayed'
$ws.Range("G105").Value = 'Dr. Amira Ibrahim, Dr. Basma Hamed'
$ws.Range("G113").Value = 'Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida'
$ws.Range("G115").Value = 'Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Yasmin, Dr. Remon, Dr. Monica, Dr. Marina Atef, Dr. Nardine'
